# Swap the "category-code" (column F) and "category-name" (column G) columns
# for every row in the worksheet, including the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$rangeF = "F1:F" + $lastRow
$rangeG = "G1:G" + $lastRow
$rangeHelper = "I1:I" + $lastRow

# Use a helper column (I) so the swap round-trips through Copy/PasteSpecial,
# which preserves the original cell types (shared-string text) instead of
# coercing numeric-looking strings into numbers.
$ws.Range($rangeF).Copy()
$ws.Range("I1").PasteSpecial(-4163)

$ws.Range($rangeG).Copy()
$ws.Range("F1").PasteSpecial(-4163)

$ws.Range($rangeHelper).Copy()
$ws.Range("G1").PasteSpecial(-4163)

$ws.Range($rangeHelper).Clear()
